# Update latest output (run 50)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: totals recalculated for the run ---
$wsSchedule.Range("E2").Value = 822.923712
$wsSchedule.Range("F2").Value = 13.60654285714286

# --- Detailed sheet: updated price/type values for run 50 ---
$wsDetailed.Range("B5").Value = 57.06003
$wsDetailed.Range("B6").Value = 40.54
$wsDetailed.Range("C7").Value = "historical"
$wsDetailed.Range("B8").Value = 57.06003
$wsDetailed.Range("C8").Value = "historical"
$wsDetailed.Range("B9").Value = 57.06003
$wsDetailed.Range("B10").Value = 57.06003
$wsDetailed.Range("B12").Value = 65
$wsDetailed.Range("B14").Value = 65
$wsDetailed.Range("B15").Value = 56.98
$wsDetailed.Range("B17").Value = 36.06029
$wsDetailed.Range("B18").Value = 16.50213
$wsDetailed.Range("B19").Value = 26.94954
$wsDetailed.Range("B20").Value = 36.06011
$wsDetailed.Range("B21").Value = 35.88
$wsDetailed.Range("B22").Value = 36.06011
$wsDetailed.Range("B24").Value = 36.06046
$wsDetailed.Range("B26").Value = 46.49741
$wsDetailed.Range("B29").Value = 36.06
$wsDetailed.Range("B30").Value = 29.64647
$wsDetailed.Range("B35").Value = 27.14114
$wsDetailed.Range("B36").Value = -0.32145
$wsDetailed.Range("B37").Value = -3.016
$wsDetailed.Range("B38").Value = -2.85235
$wsDetailed.Range("B39").Value = -2.81987
$wsDetailed.Range("B40").Value = 4.5258
$wsDetailed.Range("B41").Value = 9.627840000000001
$wsDetailed.Range("B42").Value = 29.59904
$wsDetailed.Range("B43").Value = 9.659560000000001
$wsDetailed.Range("B44").Value = 9.36539
$wsDetailed.Range("B45").Value = 9.754899999999999
$wsDetailed.Range("B47").Value = 57.06004
